$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Fill panels data completed successfully."
# The last question row (19) had its wording/answers corrected:
#  - question text changed from "increasing order" to "correct order" phrasing
#  - the sample numbers and correct-order answer key were updated
$ws.Range("A19").Value = "Βάλτε στη σωστή σειρά τους παρακάτω αριθμούς."
$ws.Range("C19").Value = "1|4|2|7|10|8|5"
$ws.Range("D19").Value = "1|2|4|5|7|8|10"

# Reflect the new scroll/selection state that was saved with the workbook
$ws.Activate()
$ws.Range("A23").Select()
